{"js": "// Applies the \"Working on web application\" edits to SkriptVideopitch.docx\n// Each edit searches for the exact original text and replaces it with the\n// new text, which keeps run-level formatting (italic etc.) intact because\n// the match lies fully inside runs that already carry that formatting.\n\nasync function replaceOnce(context, searchText, newText, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = context.document.body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1. \" vor dem Computer \u2013 Over the Shoulder oder frontal mit Bildschirm\"\n//    -> \" vor dem Computer \u2013frontal mit Bildschirm\"\nawait replaceOnce(\n  context,\n  \" vor dem Computer \u2013 Over the Shoulder oder frontal mit Bildschirm\",\n  \" vor dem Computer \u2013frontal mit Bildschirm\"\n);\n\n// 2. \"Setting: Timos Zimmer? Je nach Hintergrundunsch\u00e4rfe. Schreibtisch ein\n//     wenig vor schieben und von der Wand aus Filmen\" -> \"Setting: Timos Zimmer\"\nawait replaceOnce(\n  context,\n  \"Setting: Timos Zimmer? Je nach Hintergrundunsch\u00e4rfe. Schreibtisch ein wenig vor schieben und von der Wand aus Filmen\",\n  \"Setting: Timos Zimmer\"\n);\n\n// 3. \"Setting: Lukas\u2018 Arbeitsplatz oder B\u00fcro?\" -> \"Setting: Lukas\u2018 Arbeitsplatz oder Timos Zimmer?\"\nawait replaceOnce(\n  context,\n  \"Setting: Lukas\u2018 Arbeitsplatz oder B\u00fcro?\",\n  \"Setting: Lukas\u2018 Arbeitsplatz oder Timos Zimmer?\"\n);\n\n// 4. \"Setting: Wohnzimmerwand? B\u00fcroregale unscharf?\" -> \"Setting: Vordach Timos Schule, Notfall Waschk\u00fcche\"\nawait replaceOnce(\n  context,\n  \"Setting: Wohnzimmerwand? B\u00fcroregale unscharf?\",\n  \"Setting: Vordach Timos Schule, Notfall Waschk\u00fcche\"\n);\n\n// 5. \"Material: Kamera, Stativ, Mikrofon\" -> \"Material: Kamera, Stativ, Mikrofon, Softboxen?\"\nawait replaceOnce(\n  context,\n  \"Material: Kamera, Stativ, Mikrofon\",\n  \"Material: Kamera, Stativ, Mikrofon, Softboxen?\"\n);\n\n// 6. \"Setting: Sofa vor Wand oder Esstisch vor Fenster\" -> \"Setting: Vor Timos Fenster\"\nawait replaceOnce(\n  context,\n  \"Setting: Sofa vor Wand oder Esstisch vor Fenster\",\n  \"Setting: Vor Timos Fenster\"\n);\n\n// 7. \"Setting: Wohnzimmer mit Blick auf T\u00fcr und Lukas\u2018 Schreibtisch?\"\n//    -> \"Setting: Timos Zimmer mit Blick auf T\u00fcr und Schreibtisch\"\nawait replaceOnce(\n  context,\n  \"Setting: Wohnzimmer mit Blick auf T\u00fcr und Lukas\u2018 Schreibtisch?\",\n  \"Setting: Timos Zimmer mit Blick auf T\u00fcr und Schreibtisch\"\n);\n", "ps1": "# Applies the \"Working on web application\" edits to SkriptVideopitch.docx\n# using Word COM Find/Replace. Each call searches the whole document body\n# for the exact original phrase and substitutes the new text, which keeps\n# the surrounding run formatting (italic, etc.) intact.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1. \" vor dem Computer \u2013 Over the Shoulder oder frontal mit Bildschirm\"\n#    -> \" vor dem Computer \u2013frontal mit Bildschirm\"\nReplace-Text \" vor dem Computer \u2013 Over the Shoulder oder frontal mit Bildschirm\" \" vor dem Computer \u2013frontal mit Bildschirm\"\n\n# 2. \"Setting: Timos Zimmer? Je nach Hintergrundunsch\u00e4rfe. Schreibtisch ein\n#     wenig vor schieben und von der Wand aus Filmen\" -> \"Setting: Timos Zimmer\"\nReplace-Text \"Setting: Timos Zimmer? Je nach Hintergrundunsch\u00e4rfe. Schreibtisch ein wenig vor schieben und von der Wand aus Filmen\" \"Setting: Timos Zimmer\"\n\n# 3. \"Setting: Lukas\u2018 Arbeitsplatz oder B\u00fcro?\" -> \"Setting: Lukas\u2018 Arbeitsplatz oder Timos Zimmer?\"\nReplace-Text \"Setting: Lukas\u2018 Arbeitsplatz oder B\u00fcro?\" \"Setting: Lukas\u2018 Arbeitsplatz oder Timos Zimmer?\"\n\n# 4. \"Setting: Wohnzimmerwand? B\u00fcroregale unscharf?\" -> \"Setting: Vordach Timos Schule, Notfall Waschk\u00fcche\"\nReplace-Text \"Setting: Wohnzimmerwand? B\u00fcroregale unscharf?\" \"Setting: Vordach Timos Schule, Notfall Waschk\u00fcche\"\n\n# 5. \"Material: Kamera, Stativ, Mikrofon\" -> \"Material: Kamera, Stativ, Mikrofon, Softboxen?\"\nReplace-Text \"Material: Kamera, Stativ, Mikrofon\" \"Material: Kamera, Stativ, Mikrofon, Softboxen?\"\n\n# 6. \"Setting: Sofa vor Wand oder Esstisch vor Fenster\" -> \"Setting: Vor Timos Fenster\"\nReplace-Text \"Setting: Sofa vor Wand oder Esstisch vor Fenster\" \"Setting: Vor Timos Fenster\"\n\n# 7. \"Setting: Wohnzimmer mit Blick auf T\u00fcr und Lukas\u2018 Schreibtisch?\"\n#    -> \"Setting: Timos Zimmer mit Blick auf T\u00fcr und Schreibtisch\"\nReplace-Text \"Setting: Wohnzimmer mit Blick auf T\u00fcr und Lukas\u2018 Schreibtisch?\" \"Setting: Timos Zimmer mit Blick auf T\u00fcr und Schreibtisch\"\n"}
